$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the list date (A1)
$ws.Range("A1").Value = 45436

# Update prices in column D for rows 28-34
$ws.Range("D28").Value = 700.646
$ws.Range("D29").Value = 734.098
$ws.Range("D30").Value = 1182.716
$ws.Range("D31").Value = 1679.442
$ws.Range("D32").Value = 2346.599
$ws.Range("D33").Value = 3011.67
$ws.Range("D34").Value = 4705.737
